$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 227.8
$ws.Range("I2").Value = 227.8
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 227.8
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -114.8
$ws.Range("N2").ClearContents()
$ws.Range("H11").Value = 138.33333
$ws.Range("I11").Value = 138.33333
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 138.33333
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 1.666670000000011
$ws.Range("H18").Value = 2249.9375
$ws.Range("I18").Value = 2249.9375
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2249.9375
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1965.9375
$ws.Range("H19").Value = 1115.6923
$ws.Range("I19").Value = 775.44446
$ws.Range("J19").Value = 1295.8235
$ws.Range("K19").Value = 775.44446
$ws.Range("L19").Value = 1295.8235
$ws.Range("M19").Value = -600.44446
$ws.Range("N19").Value = -1645.8235
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H33").Value = 129.06061
$ws.Range("I33").Value = 130.28125
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 130.28125
$ws.Range("L33").Value = 90
$ws.Range("M33").Value = 98.71875
$ws.Range("N33").Value = -548
$ws.Range("H51").Value = 7119.8945
$ws.Range("I51").Value = 3966.5
$ws.Range("J51").Value = 8575.308000000001
$ws.Range("K51").Value = 3966.5
$ws.Range("L51").Value = 8575.308000000001
$ws.Range("M51").Value = -3482.5
$ws.Range("N51").Value = -9543.308000000001
$ws.Range("H58").Value = 45
$ws.Range("I58").Value = 45
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 135
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 15
$ws.Range("H88").Value = 2447.842
$ws.Range("I88").Value = 2319.5
$ws.Range("J88").Value = 2541.182
$ws.Range("K88").Value = 2319.5
$ws.Range("L88").Value = 2541.182
$ws.Range("M88").Value = -1913.5
$ws.Range("N88").Value = -3353.182
$ws.Range("H91").Value = 2447.842
$ws.Range("I91").Value = 2319.5
$ws.Range("J91").Value = 2541.182
$ws.Range("K91").Value = 2319.5
$ws.Range("L91").Value = 2541.182
$ws.Range("M91").Value = -915.5
$ws.Range("N91").Value = -5349.182
$ws.Range("H98").Value = 10499.5
$ws.Range("I98").Value = 10499.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 10499.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -9001.5
$ws.Range("H122").Value = 10499.5
$ws.Range("I122").Value = 10499.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 31498.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -29048.5
$ws.Range("H132").Value = 1435.1786
$ws.Range("I132").Value = 1345.238
$ws.Range("J132").Value = 1705
$ws.Range("K132").Value = 4035.714
$ws.Range("L132").Value = 5115
$ws.Range("M132").Value = -1505.714
$ws.Range("N132").Value = -10175
$ws.Range("H137").Value = 2026.1765
$ws.Range("I137").Value = 1443.5
$ws.Range("J137").Value = 2858.5715
$ws.Range("K137").Value = 4330.5
$ws.Range("L137").Value = 8575.7145
$ws.Range("M137").Value = -1780.5
$ws.Range("H141").Value = 4104.3
$ws.Range("I141").Value = 4338.6665
$ws.Range("J141").Value = 1995
$ws.Range("K141").Value = 13015.9995
$ws.Range("L141").Value = 5985
$ws.Range("M141").Value = -7835.999500000002
$ws.Range("N141").Value = -16345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4489.5454
$ws.Range("I2").Value = 3141.4285
$ws.Range("J2").Value = 6848.75
$ws.Range("K2").Value = 3141.4285
$ws.Range("L2").Value = 6848.75
$ws.Range("M2").Value = -3028.4285
$ws.Range("N2").Value = -7074.75
$ws.Range("H32").Value = 5367.25
$ws.Range("I32").Value = 5389.019
$ws.Range("J32").Value = 5084.25
$ws.Range("K32").Value = 5389.019
$ws.Range("L32").Value = 5084.25
$ws.Range("M32").Value = -5102.019
$ws.Range("H45").Value = 95876.67999999999
$ws.Range("I45").Value = 158108.61
$ws.Range("J45").Value = 5986.1113
$ws.Range("K45").Value = 158108.61
$ws.Range("L45").Value = 5986.1113
$ws.Range("M45").Value = -157731.61
$ws.Range("H61").Value = 14711678
$ws.Range("I61").Value = 16671502
$ws.Range("J61").Value = 13000
$ws.Range("K61").Value = 16671502
$ws.Range("L61").Value = 13000
$ws.Range("M61").Value = -16671290
$ws.Range("H110").Value = 8469.134
$ws.Range("I110").Value = 7276.1816
$ws.Range("J110").Value = 11749.75
$ws.Range("K110").Value = 7276.1816
$ws.Range("L110").Value = 11749.75
$ws.Range("M110").Value = -5231.1816
$ws.Range("H116").Value = 4489.5454
$ws.Range("I116").Value = 3141.4285
$ws.Range("J116").Value = 6848.75
$ws.Range("K116").Value = 3141.4285
$ws.Range("L116").Value = 6848.75
$ws.Range("M116").Value = -847.4285
$ws.Range("N116").Value = -11436.75
$ws.Range("H132").Value = 4442.815
$ws.Range("I132").Value = 3940.2917
$ws.Range("J132").Value = 8463
$ws.Range("K132").Value = 11820.8751
$ws.Range("L132").Value = 25389
$ws.Range("M132").Value = -9290.875100000001
$ws.Range("H136").Value = 14711678
$ws.Range("I136").Value = 16671502
$ws.Range("J136").Value = 13000
$ws.Range("K136").Value = 50014506
$ws.Range("L136").Value = 39000
$ws.Range("M136").Value = -50011956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4489.5454
$ws.Range("I3").Value = 3141.4285
$ws.Range("J3").Value = 6848.75
$ws.Range("K3").Value = 3141.4285
$ws.Range("L3").Value = 6848.75
$ws.Range("M3").Value = -3027.4285
$ws.Range("N3").Value = -7076.75
$ws.Range("H46").Value = 30001
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 30001
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 30001
$ws.Range("N46").Value = -30597
$ws.Range("H80").Value = 691.2222
$ws.Range("I80").Value = 488.66666
$ws.Range("J80").Value = 792.5
$ws.Range("K80").Value = 488.66666
$ws.Range("L80").Value = 792.5
$ws.Range("M80").Value = 509.33334
$ws.Range("N80").Value = -2788.5
$ws.Range("H83").Value = 691.2222
$ws.Range("I83").Value = 488.66666
$ws.Range("J83").Value = 792.5
$ws.Range("K83").Value = 2443.3333
$ws.Range("L83").Value = 3962.5
$ws.Range("M83").Value = 2548.6667
$ws.Range("N83").Value = -13946.5
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 15584.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 15584.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 15584.75
$ws.Range("N100").Value = -17748.75
$ws.Range("H107").Value = 6360
$ws.Range("I107").Value = 5004.4
$ws.Range("J107").Value = 9749
$ws.Range("K107").Value = 5004.4
$ws.Range("L107").Value = 9749
$ws.Range("M107").Value = -3084.4
$ws.Range("N107").Value = -13589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3808.5144
$ws.Range("I31").Value = 2730.25
$ws.Range("J31").Value = 6161.091
$ws.Range("K31").Value = 2730.25
$ws.Range("L31").Value = 6161.091
$ws.Range("M31").Value = -2435.25
$ws.Range("N31").Value = -6751.091
$ws.Range("H34").Value = 3808.5144
$ws.Range("I34").Value = 2730.25
$ws.Range("J34").Value = 6161.091
$ws.Range("K34").Value = 2730.25
$ws.Range("L34").Value = 6161.091
$ws.Range("M34").Value = -2528.25
$ws.Range("N34").Value = -6565.091
$ws.Range("H92").Value = 50000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 50000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H96").Value = 50000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 50000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H132").Value = 6608.222
$ws.Range("I132").Value = 6218.5713
$ws.Range("J132").Value = 7972
$ws.Range("K132").Value = 18655.7139
$ws.Range("L132").Value = 23916
$ws.Range("M132").Value = -16125.7139
$ws.Range("N132").Value = -28976

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5140.1875
$ws.Range("I3").Value = 5584.8
$ws.Range("J3").Value = 4399.1665
$ws.Range("K3").Value = 16754.4
$ws.Range("L3").Value = 13197.4995
$ws.Range("M3").Value = -16642.4
$ws.Range("H11").Value = 9091325
$ws.Range("I11").Value = 10000357
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 30001071
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -30000931
$ws.Range("H12").Value = 1073
$ws.Range("I12").Value = 635.75
$ws.Range("J12").Value = 1267.3334
$ws.Range("K12").Value = 1907.25
$ws.Range("L12").Value = 3802.0002
$ws.Range("M12").Value = -1734.25
$ws.Range("N12").Value = -4148.0002
$ws.Range("H60").Value = 395
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 395
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 1185
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -1687
$ws.Range("H81").Value = 7610.7
$ws.Range("I81").Value = 937.6667
$ws.Range("J81").Value = 10470.571
$ws.Range("K81").Value = 2813.0001
$ws.Range("L81").Value = 31411.713
$ws.Range("M81").Value = -1690.0001
$ws.Range("N81").Value = -33657.713
$ws.Range("H84").Value = 7610.7
$ws.Range("I84").Value = 937.6667
$ws.Range("J84").Value = 10470.571
$ws.Range("K84").Value = 8439.0003
$ws.Range("L84").Value = 94235.139
$ws.Range("M84").Value = -2823.0003
$ws.Range("N84").Value = -105467.139
$ws.Range("H130").Value = 10666.667
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 15000
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 45000
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -55040
$ws.Range("H137").Value = 6371.5
$ws.Range("I137").Value = 1503.1818
$ws.Range("J137").Value = 24222
$ws.Range("K137").Value = 4509.5454
$ws.Range("L137").Value = 72666
$ws.Range("M137").Value = 590.4546
$ws.Range("H139").Value = 1885.3043
$ws.Range("I139").Value = 1885.3043
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5655.9129
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -515.9129000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 4883.8066
$ws.Range("I132").Value = 3488.6287
$ws.Range("J132").Value = 6692.3706
$ws.Range("K132").Value = 10465.8861
$ws.Range("L132").Value = 20077.1118
$ws.Range("M132").Value = -7935.8861
$ws.Range("N132").Value = -25137.1118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 40983
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 40983
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 40983
$ws.Range("N104").Value = -47971
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 6575.1665
$ws.Range("I122").Value = 3968.6667
$ws.Range("J122").Value = 7444
$ws.Range("K122").Value = 11906.0001
$ws.Range("L122").Value = 22332
$ws.Range("M122").Value = -9456.000100000001
$ws.Range("N122").Value = -27232
$ws.Range("H132").Value = 2964.3513
$ws.Range("I132").Value = 2596.2424
$ws.Range("J132").Value = 6001.25
$ws.Range("K132").Value = 7788.7272
$ws.Range("L132").Value = 18003.75
$ws.Range("M132").Value = -5258.7272
$ws.Range("N132").Value = -23063.75
